# Update crypto price (D) and 1h volume change (E) columns for rows 2-51
# with freshly scraped values from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.847.34"
$ws.Range("E2").Value = "  +5.04%  "
$ws.Range("D3").Value = "1.894.57"
$ws.Range("E3").Value = "  +3.87%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("E4").Value = "  -0.69%  "
$ws.Range("D5").Value = "'339.13"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").Value = "'0.9989"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "'0.4735"
$ws.Range("E7").Value = "  +3.16%  "
$ws.Range("D8").Value = "'0.4040"
$ws.Range("E8").Value = "  +6.07%  "
$ws.Range("D9").Value = "'47.71"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").Value = "'0.08102"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("D11").Value = "'1.014"
$ws.Range("E11").Value = "  +4.74%  "
$ws.Range("D12").Value = "'22.30"
$ws.Range("E12").Value = "  +6.12%  "
$ws.Range("D13").Value = "'6.077"
$ws.Range("E13").Value = "  +3.29%  "
$ws.Range("D14").Value = "1.882.98"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").Value = "'7.334"
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("D16").Value = "'90.97"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "'0.9997"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'0.00001052"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").Value = "'0.06601"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'17.71"
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("D21").Value = "'0.9986"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "28.867.93"
$ws.Range("E22").Value = "  +5.20%  "
$ws.Range("D23").Value = "'5.530"
$ws.Range("E23").Value = "  +3.70%  "
$ws.Range("D24").Value = "'11.10"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "'2.263"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "2.120.96"
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("D27").Value = "'160.56"
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").Value = "'19.93"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("D29").Value = "'2.149"
$ws.Range("E29").Value = "  +4.55%  "
$ws.Range("D30").Value = "'5.538"
$ws.Range("E30").Value = "  +5.08%  "
$ws.Range("D31").Value = "'120.32"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "'1.004"
$ws.Range("E32").Value = "  +6.37%  "
$ws.Range("D33").Value = "'0.09574"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").Value = "'3.649"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").Value = "'1.400"
$ws.Range("E35").Value = "  +6.17%  "
$ws.Range("D36").Value = "'5.395"
$ws.Range("E36").Value = "  +3.01%  "
$ws.Range("D37").Value = "'0.06191"
$ws.Range("E37").Value = "  +4.54%  "
$ws.Range("D38").Value = "'0.02286"
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("D39").Value = "'8.574"
$ws.Range("E39").Value = "  +6.84%  "
$ws.Range("D40").Value = "'1.191"
$ws.Range("E40").Value = "  +2.62%  "
$ws.Range("D41").Value = "'0.5998"
$ws.Range("E41").Value = "  +4.33%  "
$ws.Range("D42").Value = "'0.1898"
$ws.Range("E42").Value = "  +3.89%  "
$ws.Range("D43").Value = "'0.9982"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "'10.40"
$ws.Range("E44").Value = "  +3.84%  "
$ws.Range("D45").Value = "'1.264"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "'0.5610"
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("D47").Value = "'12.29"
$ws.Range("E47").Value = "  +2.94%  "
$ws.Range("D48").Value = "'1.973"
$ws.Range("E48").Value = "  +5.73%  "
$ws.Range("D49").Value = "'0.07245"
$ws.Range("E49").Value = "  +9.63%  "
$ws.Range("D50").Value = "'2.122"
$ws.Range("E50").Value = "  +15.25%  "
$ws.Range("D51").Value = "'112.81"
$ws.Range("E51").Value = "  +2.14%  "
